$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A57").Value = "NAAAR"
$ws.Range("B57").Value = "N 52:00:00.00E 03:00:00.00"

$ws.Range("B57").Select() | Out-Null
